# Revert to 2.1.1 files
# -------------------------------------------------------------------------
# This script reproduces (via Excel COM automation) the diff that:
#   1. Inserts a new "Texas Notes" worksheet between "About" and
#      "MSCdtRPbQL" with reviewer notes about the rebate-program source.
#   2. Updates the hyperlink display text on the "About" sheet (B6) to a
#      new source URL.
#   3. Updates the "MSCdtRPbQL" sheet's rebate-qualifying change value
#      from 7.4% to 4.95% (the dependent formula in B2 recalculates
#      automatically).
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("MSCdtRPbQL")

# ---------------------------------------------------------------------
# 0. "About" sheet: refresh the displayed hyperlink text (B6) first so
#    the new shared string lands near the top of the table, matching
#    the reference document's ordering.
# ---------------------------------------------------------------------
$aboutSheet.Range("B6").Value = "http://www.cepe.ethz.ch/publications/workingPapers/CEPE_WP86.pdf"

# ---------------------------------------------------------------------
# 1. Insert the new "Texas Notes" worksheet right before "MSCdtRPbQL"
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Add($dataSheet)
$notes.Name = "Texas Notes"

$notes.Range("A1").Value  = "I found an updated version of the study that EPS cites here:"
$notes.Range("A2").Value  = "DOI:10.1007/s12053-015-9386-7"
$notes.Range("A2").Font.Name  = "Arial"
$notes.Range("A2").Font.Color = 5592405

$notes.Range("A4").Value  = """"
$notes.Range("A5").Value  = "Results suggest that rebate policies increase"
$notes.Range("A6").Value  = "the sales share of ENERGY STAR household appliances"
$notes.Range("A7").Value  = "by 3.3 to 6.6 percentage points, and this represents"
$notes.Range("A8").Value  = "an impact of 9 to 18 % on the mean level of the"
$notes.Range("A9").Value  = "sales share of ENERGY STAR household appliances"
$notes.Range("A10").Value = "in the US between 2001 and 2006."
$notes.Range("A11").Value = """"
$notes.Range("A12").Value = "See Table 3."

$notes.Range("A14").Value = "The 3.3% and 6.6% answers are both statistically significant and depend on the analysis method being used. "
$notes.Range("A15").Value = "Neither method stood out to me as being ""better"", so maybe we just average the two together and get 4.95%."

$notes.Range("A17").Value = "And, for what it's worth, it looks like the 7.4% number that EPS was using comes from Table 1 and does not mean what they intended it to mean. "
$notes.Range("A18").Value = "The 7.4% describes the mean of a distribution of whether or not a state has rebates for particular applicances.  "
$notes.Range("A19").Value = "It's a binary distribution, so 0 for states without rebates and 1 for states with rebates. "

$notes.Range("A21").Value = "The 4.95% number (average of 3.3 and 6.6) has to do with the increase in energy start appliances when a rebate is available. "

$notes.Range("A24").Value = "***I looked a bit for Texas specific data and didn't find much. "
$notes.Range("A25").Value = "It's worth noting that Texas does not appear to be as rebate-happy as most states. So there is probably not a lot of historical data for a Texas-specific study"
$notes.Range("A26").Value = "to build on. Then, I assume that Texas consumers will respond to rebates in a similar way to American consumers. "

# ---------------------------------------------------------------------
# 3. "MSCdtRPbQL" sheet: update the rebate-qualifying change value
# ---------------------------------------------------------------------
$dataSheet.Range("C2").Value = 0.0495

# ---------------------------------------------------------------------
# 4. Selections / active sheet, matching the reviewer's final view state
# ---------------------------------------------------------------------
$aboutSheet.Activate()
$aboutSheet.Range("A11").Select()

$dataSheet.Activate()
$dataSheet.Range("C3").Select()

$notes.Activate()
$notes.Range("A27").Select()
